# -----------------------------------------------------------------------
# Add the new "2022-Q3" quarter sheet (data for the newest quarter) and
# place it right after the "总计" (summary) sheet, i.e. as the new 2nd tab.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$total    = $wb.Worksheets.Item("总计")
$refSheet = $wb.Worksheets.Item("2022-Q2")

# Create the sheet first (default placement), populate it while it is the
# "fresh" sheet (formatting copy/paste is only reliable on a sheet that
# hasn't been repositioned yet), then move it into its final slot.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"

# ---- header row (row 1): copy style from the 2022-Q2 sheet, then set text
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- index column (A2:A36): copy the bold/bordered index-cell style
$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A36").PasteSpecial(-4122)

# ---- fund rows (35 funds, rows 2-36)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'" + "161726"
$newSheet.Range("C2").Value = "招商国证生物医药指数A"
$newSheet.Range("D2").Value = "'" + "117.97"
$newSheet.Range("E2").Value = "'" + "94.84"
$newSheet.Range("F2").Value = "'" + "4.47"
$newSheet.Range("G2").Value = "'" + "5.2733"
$newSheet.Range("H2").Value = 7
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'" + "159992"
$newSheet.Range("C3").Value = "银华中证创新药产业ETF"
$newSheet.Range("D3").Value = "'" + "43.46"
$newSheet.Range("E3").Value = "'" + "98.40"
$newSheet.Range("F3").Value = "'" + "3.86"
$newSheet.Range("G3").Value = "'" + "1.6776"
$newSheet.Range("H3").Value = 7
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'" + "001257"
$newSheet.Range("C4").Value = "兴业收益增强债券A"
$newSheet.Range("D4").Value = "'" + "68.35"
$newSheet.Range("E4").Value = "'" + "20.08"
$newSheet.Range("F4").Value = "'" + "1.39"
$newSheet.Range("G4").Value = "'" + "0.9501"
$newSheet.Range("H4").Value = 7
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'" + "012417"
$newSheet.Range("C5").Value = "招商国证生物医药指数C"
$newSheet.Range("D5").Value = "'" + "16.74"
$newSheet.Range("E5").Value = "'" + "94.84"
$newSheet.Range("F5").Value = "'" + "4.47"
$newSheet.Range("G5").Value = "'" + "0.7483"
$newSheet.Range("H5").Value = 7
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'" + "501009"
$newSheet.Range("C6").Value = "汇添富中证生物科技主题指数（LOF）A"
$newSheet.Range("D6").Value = "'" + "17.64"
$newSheet.Range("E6").Value = "'" + "94.60"
$newSheet.Range("F6").Value = "'" + "3.92"
$newSheet.Range("G6").Value = "'" + "0.6915"
$newSheet.Range("H6").Value = 8
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'" + "159859"
$newSheet.Range("C7").Value = "天弘国证生物医药ETF"
$newSheet.Range("D7").Value = "'" + "14.46"
$newSheet.Range("E7").Value = "'" + "99.79"
$newSheet.Range("F7").Value = "'" + "4.71"
$newSheet.Range("G7").Value = "'" + "0.6811"
$newSheet.Range("H7").Value = 7
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'" + "515120"
$newSheet.Range("C8").Value = "广发中证创新药产业ETF"
$newSheet.Range("D8").Value = "'" + "16.65"
$newSheet.Range("E8").Value = "'" + "99.40"
$newSheet.Range("F8").Value = "'" + "3.90"
$newSheet.Range("G8").Value = "'" + "0.6494"
$newSheet.Range("H8").Value = 7
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'" + "159837"
$newSheet.Range("C9").Value = "易方达中证生物科技主题ETF"
$newSheet.Range("D9").Value = "'" + "16.12"
$newSheet.Range("E9").Value = "'" + "99.10"
$newSheet.Range("F9").Value = "'" + "3.93"
$newSheet.Range("G9").Value = "'" + "0.6335"
$newSheet.Range("H9").Value = 8
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'" + "501010"
$newSheet.Range("C10").Value = "汇添富中证生物科技主题指数（LOF）C"
$newSheet.Range("D10").Value = "'" + "13.43"
$newSheet.Range("E10").Value = "'" + "94.60"
$newSheet.Range("F10").Value = "'" + "3.92"
$newSheet.Range("G10").Value = "'" + "0.5265"
$newSheet.Range("H10").Value = 8
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "'" + "005984"
$newSheet.Range("C11").Value = "兴业聚华混合A"
$newSheet.Range("D11").Value = "'" + "23.61"
$newSheet.Range("E11").Value = "'" + "29.42"
$newSheet.Range("F11").Value = "'" + "1.98"
$newSheet.Range("G11").Value = "'" + "0.4675"
$newSheet.Range("H11").Value = 7
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "'" + "159839"
$newSheet.Range("C12").Value = "汇添富国证生物医药ETF"
$newSheet.Range("D12").Value = "'" + "3.99"
$newSheet.Range("E12").Value = "'" + "99.64"
$newSheet.Range("F12").Value = "'" + "4.72"
$newSheet.Range("G12").Value = "'" + "0.1883"
$newSheet.Range("H12").Value = 7
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "'" + "001258"
$newSheet.Range("C13").Value = "兴业收益增强债券C"
$newSheet.Range("D13").Value = "'" + "11.01"
$newSheet.Range("E13").Value = "'" + "20.08"
$newSheet.Range("F13").Value = "'" + "1.39"
$newSheet.Range("G13").Value = "'" + "0.1530"
$newSheet.Range("H13").Value = 7
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "'" + "005985"
$newSheet.Range("C14").Value = "兴业聚华混合C"
$newSheet.Range("D14").Value = "'" + "6.45"
$newSheet.Range("E14").Value = "'" + "29.42"
$newSheet.Range("F14").Value = "'" + "1.98"
$newSheet.Range("G14").Value = "'" + "0.1277"
$newSheet.Range("H14").Value = 7
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "'" + "517850"
$newSheet.Range("C15").Value = "汇添富中证沪港深张江自主创新50ETF"
$newSheet.Range("D15").Value = "'" + "2.26"
$newSheet.Range("E15").Value = "'" + "97.75"
$newSheet.Range("F15").Value = "'" + "5.03"
$newSheet.Range("G15").Value = "'" + "0.1137"
$newSheet.Range("H15").Value = 5
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "'" + "516080"
$newSheet.Range("C16").Value = "易方达中证创新药产业ETF"
$newSheet.Range("D16").Value = "'" + "2.47"
$newSheet.Range("E16").Value = "'" + "98.96"
$newSheet.Range("F16").Value = "'" + "3.87"
$newSheet.Range("G16").Value = "'" + "0.0956"
$newSheet.Range("H16").Value = 7
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "'" + "515960"
$newSheet.Range("C17").Value = "嘉实中证医药健康100策略ETF"
$newSheet.Range("D17").Value = "'" + "2.99"
$newSheet.Range("E17").Value = "'" + "98.86"
$newSheet.Range("F17").Value = "'" + "2.65"
$newSheet.Range("G17").Value = "'" + "0.0792"
$newSheet.Range("H17").Value = 8
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "'" + "512120"
$newSheet.Range("C18").Value = "华安中证细分医药ETF"
$newSheet.Range("D18").Value = "'" + "2.93"
$newSheet.Range("E18").Value = "'" + "98.58"
$newSheet.Range("F18").Value = "'" + "2.69"
$newSheet.Range("G18").Value = "'" + "0.0788"
$newSheet.Range("H18").Value = 10
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "'" + "517120"
$newSheet.Range("C19").Value = "华泰柏瑞中证沪港深创新药产业ETF"
$newSheet.Range("D19").Value = "'" + "2.96"
$newSheet.Range("E19").Value = "'" + "94.48"
$newSheet.Range("F19").Value = "'" + "2.52"
$newSheet.Range("G19").Value = "'" + "0.0746"
$newSheet.Range("H19").Value = 10
$newSheet.Range("A20").Value = 18
$newSheet.Range("B20").Value = "'" + "165519"
$newSheet.Range("C20").Value = "信诚中证800医药指数（LOF）A"
$newSheet.Range("D20").Value = "'" + "2.63"
$newSheet.Range("E20").Value = "'" + "93.95"
$newSheet.Range("F20").Value = "'" + "2.63"
$newSheet.Range("G20").Value = "'" + "0.0692"
$newSheet.Range("H20").Value = 10
$newSheet.Range("A21").Value = 19
$newSheet.Range("B21").Value = "'" + "560900"
$newSheet.Range("C21").Value = "上投摩根中证创新药产业ETF"
$newSheet.Range("D21").Value = "'" + "1.93"
$newSheet.Range("E21").Value = "'" + "94.08"
$newSheet.Range("F21").Value = "'" + "3.51"
$newSheet.Range("G21").Value = "'" + "0.0677"
$newSheet.Range("H21").Value = 7
$newSheet.Range("A22").Value = 20
$newSheet.Range("B22").Value = "'" + "517380"
$newSheet.Range("C22").Value = "天弘恒生沪深港创新药精选50ETF"
$newSheet.Range("D22").Value = "'" + "1.98"
$newSheet.Range("E22").Value = "'" + "99.42"
$newSheet.Range("F22").Value = "'" + "3.39"
$newSheet.Range("G22").Value = "'" + "0.0671"
$newSheet.Range("H22").Value = 8
$newSheet.Range("A23").Value = 21
$newSheet.Range("B23").Value = "'" + "159838"
$newSheet.Range("C23").Value = "博时中证医药50ETF"
$newSheet.Range("D23").Value = "'" + "1.59"
$newSheet.Range("E23").Value = "'" + "98.09"
$newSheet.Range("F23").Value = "'" + "2.67"
$newSheet.Range("G23").Value = "'" + "0.0425"
$newSheet.Range("H23").Value = 10
$newSheet.Range("A24").Value = 22
$newSheet.Range("B24").Value = "'" + "159748"
$newSheet.Range("C24").Value = "富国中证沪港深创新药产业ETF"
$newSheet.Range("D24").Value = "'" + "1.62"
$newSheet.Range("E24").Value = "'" + "98.42"
$newSheet.Range("F24").Value = "'" + "2.53"
$newSheet.Range("G24").Value = "'" + "0.0410"
$newSheet.Range("H24").Value = 10
$newSheet.Range("A25").Value = 23
$newSheet.Range("B25").Value = "'" + "159849"
$newSheet.Range("C25").Value = "招商中证生物科技主题ETF"
$newSheet.Range("D25").Value = "'" + "0.99"
$newSheet.Range("E25").Value = "'" + "99.25"
$newSheet.Range("F25").Value = "'" + "3.98"
$newSheet.Range("G25").Value = "'" + "0.0394"
$newSheet.Range("H25").Value = 8
$newSheet.Range("A26").Value = 24
$newSheet.Range("B26").Value = "'" + "510660"
$newSheet.Range("C26").Value = "华夏上证医药卫生ETF"
$newSheet.Range("D26").Value = "'" + "0.96"
$newSheet.Range("E26").Value = "'" + "99.53"
$newSheet.Range("F26").Value = "'" + "4.08"
$newSheet.Range("G26").Value = "'" + "0.0392"
$newSheet.Range("H26").Value = 4
$newSheet.Range("A27").Value = 25
$newSheet.Range("B27").Value = "'" + "516060"
$newSheet.Range("C27").Value = "工银新药ETF"
$newSheet.Range("D27").Value = "'" + "0.99"
$newSheet.Range("E27").Value = "'" + "98.53"
$newSheet.Range("F27").Value = "'" + "3.84"
$newSheet.Range("G27").Value = "'" + "0.0380"
$newSheet.Range("H27").Value = 7
$newSheet.Range("A28").Value = 26
$newSheet.Range("B28").Value = "'" + "516500"
$newSheet.Range("C28").Value = "华夏中证生物科技主题ETF"
$newSheet.Range("D28").Value = "'" + "0.88"
$newSheet.Range("E28").Value = "'" + "98.95"
$newSheet.Range("F28").Value = "'" + "3.93"
$newSheet.Range("G28").Value = "'" + "0.0346"
$newSheet.Range("H28").Value = 8
$newSheet.Range("A29").Value = 27
$newSheet.Range("B29").Value = "'" + "159858"
$newSheet.Range("C29").Value = "南方中证创新药产业ETF"
$newSheet.Range("D29").Value = "'" + "0.76"
$newSheet.Range("E29").Value = "'" + "99.70"
$newSheet.Range("F29").Value = "'" + "3.91"
$newSheet.Range("G29").Value = "'" + "0.0297"
$newSheet.Range("H29").Value = 7
$newSheet.Range("A30").Value = 28
$newSheet.Range("B30").Value = "'" + "159835"
$newSheet.Range("C30").Value = "建信中证创新药产业ETF"
$newSheet.Range("D30").Value = "'" + "0.59"
$newSheet.Range("E30").Value = "'" + "98.09"
$newSheet.Range("F30").Value = "'" + "3.78"
$newSheet.Range("G30").Value = "'" + "0.0223"
$newSheet.Range("H30").Value = 7
$newSheet.Range("A31").Value = 29
$newSheet.Range("B31").Value = "'" + "516930"
$newSheet.Range("C31").Value = "民生加银中证生物科技主题ETF"
$newSheet.Range("D31").Value = "'" + "0.57"
$newSheet.Range("E31").Value = "'" + "97.99"
$newSheet.Range("F31").Value = "'" + "3.71"
$newSheet.Range("G31").Value = "'" + "0.0211"
$newSheet.Range("H31").Value = 8
$newSheet.Range("A32").Value = 30
$newSheet.Range("B32").Value = "'" + "012599"
$newSheet.Range("C32").Value = "华安国证生物医药指数C"
$newSheet.Range("D32").Value = "'" + "0.22"
$newSheet.Range("E32").Value = "'" + "94.33"
$newSheet.Range("F32").Value = "'" + "4.51"
$newSheet.Range("G32").Value = "'" + "0.0099"
$newSheet.Range("H32").Value = 7
$newSheet.Range("A33").Value = 31
$newSheet.Range("B33").Value = "'" + "014129"
$newSheet.Range("C33").Value = "西藏东财中证沪港深创新药产业指数C"
$newSheet.Range("D33").Value = "'" + "0.31"
$newSheet.Range("E33").Value = "'" + "94.83"
$newSheet.Range("F33").Value = "'" + "2.57"
$newSheet.Range("G33").Value = "'" + "0.0080"
$newSheet.Range("H33").Value = 10
$newSheet.Range("A34").Value = 32
$newSheet.Range("B34").Value = "'" + "012598"
$newSheet.Range("C34").Value = "华安国证生物医药指数A"
$newSheet.Range("D34").Value = "'" + "0.16"
$newSheet.Range("E34").Value = "'" + "94.33"
$newSheet.Range("F34").Value = "'" + "4.51"
$newSheet.Range("G34").Value = "'" + "0.0072"
$newSheet.Range("H34").Value = 7
$newSheet.Range("A35").Value = 33
$newSheet.Range("B35").Value = "'" + "014128"
$newSheet.Range("C35").Value = "西藏东财中证沪港深创新药产业指数A"
$newSheet.Range("D35").Value = "'" + "0.24"
$newSheet.Range("E35").Value = "'" + "94.83"
$newSheet.Range("F35").Value = "'" + "2.57"
$newSheet.Range("G35").Value = "'" + "0.0062"
$newSheet.Range("H35").Value = 10
$newSheet.Range("A36").Value = 34
$newSheet.Range("B36").Value = "'" + "013080"
$newSheet.Range("C36").Value = "信诚中证800医药指数（LOF）C"
$newSheet.Range("D36").Value = "'" + "0.12"
$newSheet.Range("E36").Value = "'" + "93.95"
$newSheet.Range("F36").Value = "'" + "2.63"
$newSheet.Range("G36").Value = "'" + "0.0032"
$newSheet.Range("H36").Value = 10

# ---- move the populated sheet into its final tab position
# (re-fetch "总计" here: a reference captured before Worksheets.Add() can
#  go stale once the sheet collection changes, which silently turns the
#  Move below into a no-op)
$totalForMove = $wb.Worksheets.Item("总计")
$newSheet.Move([System.Reflection.Missing]::Value, $totalForMove)

# -----------------------------------------------------------------------
# Update the "总计" (summary) sheet: add a 2022-Q3 entry at the top and
# push the rest of the history down by one row. The leftmost "A" index
# column is a plain running counter (0,1,2,...) tied to the row number,
# not to the quarter, so a brand-new row 9 is created (copying A8's
# style for the new A9 index cell) while columns B:D simply shift down.
# -----------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 35
$total.Range("D2").Value = 13.76

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 52
$total.Range("D3").Value = 20.3

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 70
$total.Range("D4").Value = 24.8

$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 32
$total.Range("D5").Value = 16.44

$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 33
$total.Range("D6").Value = 13.02

$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 89
$total.Range("D7").Value = 40.44

$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 17
$total.Range("D8").Value = 10.47

$total.Range("A9").Value = 7
$total.Range("B9").Value = "2020-Q4"
$total.Range("C9").Value = 26
$total.Range("D9").Value = 16.65
